$d = $word.ActiveDocument
$r = $d.Paragraphs.Item(100).Range
Write-Host "Text=[$($r.Text)]"
Write-Host "Font=$($r.Font.Color)"
